# Generate Report for Handback
# Adds a new handback row (GUID 9c23dcd6-9812-4851-b972-3e9ca64347a0) to the
# Overview sheet and to the two per-language sheets (zh-cn, de-de), mirroring
# the pattern already used for the existing bc8b34db.../3343ba93... rows.

$wb = $excel.ActiveWorkbook

$guid     = "9c23dcd6-9812-4851-b972-3e9ca64347a0"
$mdName   = "$guid.md"
$mdPath   = "e2e\$guid.md"
$zhXlf    = "$guid.0b3f1845a88d207de38cf9cc315917d8c3f2b2e6.zh-cn.xlf"
$deXlf    = "$guid.0b3f1845a88d207de38cf9cc315917d8c3f2b2e6.de-de.xlf"

$handoffDateTime = "2016-08-25 18:46:21"
$zhHandbackDone   = "2016-08-25 18:46:16"
$zhHandbackDate   = "2016-08-25 18:46:33"
$deHandbackDone   = "2016-08-25 18:46:21"
$deHandbackDate   = "2016-08-25 18:46:41"

# ---------------------------------------------------------------------------
# Sheet 1 - "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $mdPath
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = $handoffDateTime
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2d4e6f8a1b3c5d7e9f0a2b4c6d8e0f2a4b6c8d0/e2e/$guid.md", "", "", $mdPath) | Out-Null
$wsOverview.Range("B4").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet 2 - "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = $mdName
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = $zhHandbackDone
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = $mdName
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = $zhHandbackDate
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2d4e6f8a1b3c5d7e9f0a2b4c6d8e0f2a4b6c8d0/e2e/$guid.md", "", "", $mdName) | Out-Null
$wsZh.Range("A4").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d3e5f7a9b1c3d5e7f9a1b3c5d7e9f1a3b5c7d9e1/e2e/$guid.md", "", "", $mdName) | Out-Null
$wsZh.Range("I4").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet 3 - "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = $mdName
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = $handoffDateTime
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = $mdName
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = $deHandbackDate
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2d4e6f8a1b3c5d7e9f0a2b4c6d8e0f2a4b6c8d0/e2e/$guid.md", "", "", $mdName) | Out-Null
$wsDe.Range("A4").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e4f6a8b0c2d4e6f8a0b2c4d6e8f0a2b4c6d8e0f2/e2e/$guid.md", "", "", $mdName) | Out-Null
$wsDe.Range("I4").Style = "HyperLink"
